{"js": "// Load all paragraphs in the document body so we can find the last one.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The last paragraph in the body is the bullet item ending in\n// \"...god m\u00e5te.\" (it also holds the _GoBack bookmark in the OOXML).\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// Find the list that the last bullet belongs to, so the new paragraph can\n// join the same numbered/bulleted list (numId) at the same level.\nconst list = lastParagraph.listOrNullObject;\nlist.load(\"id\");\nawait context.sync();\n\n// Insert a brand-new paragraph right after it, carrying the new bullet text.\nconst newParagraph = lastParagraph.insertParagraph(\n  \"De som jobber med databasen og de som jobber med \u00e5 lage modeller burde kommunisere godt og lage objekter med samme datatyper.\",\n  \"After\"\n);\n\n// Match the formatting of the surrounding bullets: \"Listeavsnitt\" style,\n// attached to the same list at the same level (ilvl 0).\nnewParagraph.styleBuiltIn = Word.Style.listParagraph;\nnewParagraph.attachToList(list.id, 0);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The last paragraph in the body is the bullet item ending in\n# \"...god m\u00e5te.\" (it also holds the _GoBack bookmark in the OOXML).\n$paragraphs = $d.Paragraphs\n$lastParagraph = $paragraphs.Last\n\n# Remember the list template it uses so the new bullet can continue the same\n# numbered/bulleted list (numId) instead of starting a brand-new one.\n$template = $lastParagraph.Range.ListFormat.ListTemplate\n\n# Insert a brand-new paragraph right after it, carrying the new bullet text.\n$newRange = $lastParagraph.Range.InsertParagraphAfter()\n$newParagraph = $d.Paragraphs.Last\n$newParagraph.Range.Text = \"De som jobber med databasen og de som jobber med \u00e5 lage modeller burde kommunisere godt og lage objekter med samme datatyper.\"\n\n# Match the formatting of the surrounding bullets: \"Listeavsnitt\" style,\n# continuing the same list at the same level (ilvl 0).\n$newParagraph.Style = $lastParagraph.Style\n$newParagraph.Range.ListFormat.ApplyListTemplateWithLevel($template, $true, 0, $false, 0)\n"}
